$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(80).Insert()

$ws.Cells.Item(80, 1).Value = 8
$ws.Cells.Item(80, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(80, 3).Value = "Coquimbo"
$ws.Cells.Item(80, 4).Value = 44586
$ws.Cells.Item(80, 5).Value = 4
$ws.Cells.Item(80, 6).Value = 100112031
$ws.Cells.Item(80, 7).Value = "Poroto verde"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 600
$ws.Cells.Item(80, 11).Value = 30000
$ws.Cells.Item(80, 12).Value = 31000
$ws.Cells.Item(80, 13).Value = 30500
$ws.Cells.Item(80, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(80, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(80, 16).Value = 1220
$ws.Cells.Item(80, 17).Value = 25
$ws.Cells.Item(80, 18).Value = "Hortaliza"
